# Fault injection for parity generators:
# - Remove the "FAULT INJECTION LOCATION" column (Q). This shifts the
#   following columns ("EVEN ODD" and "NOTE") one position to the left,
#   so the sheet's used range shrinks from A1:S5 to A1:R5.
# - Mark fault injection as active ("YES") for every data row in the
#   (now shifted) "FAULT INJECTION" column (P).
# - Re-enter the bit-width values in columns N and O as real numbers
#   instead of text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Deleting the entire column shifts everything to its right one column
# to the left, which is exactly what turns Q/R/S (old) into Q/R (new).
$ws.Range("Q1").EntireColumn.Delete()

$lastRow = 5
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 16).Value = "YES"   # column P: FAULT INJECTION
}

$bitWidths = @{
    2 = @(36, 36)
    3 = @(36, 36)
    4 = @(256, 128)
    5 = @(256, 128)
}

foreach ($r in $bitWidths.Keys) {
    $pair = $bitWidths[$r]
    $ws.Cells.Item($r, 14).Value = $pair[0]   # column N: BIT WIDTH
    $ws.Cells.Item($r, 15).Value = $pair[1]   # column O: PARITY SOURCE BIT WIDTH
}
